$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is a simple table: col A = code, col B = status,
# col C = codeforiati:group-name, col D = codeforiati:group-code.
# The edit swaps columns C and D (both the header and every data row)
# so the group-code comes before the group-name.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 1; $r -le $lastRow; $r++) {
    $cCell = $ws.Cells.Item($r, 3)
    $dCell = $ws.Cells.Item($r, 4)
    $cVal = $cCell.Value2
    $dVal = $dCell.Value2
    $cCell.Value2 = $dVal
    $dCell.Value2 = $cVal
}
